$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ToDO")

# Row 15 gains a new "to do" entry (mirrors the layout/format of row 14).
# Copy D14's number format onto D15 first so the new date cell reuses the
# existing date style instead of minting a new one.
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A15").Value = "Connectieklasse databank gemaakt (.java)"
$ws.Range("B15").Value = "1 uur"
$ws.Range("C15").Value = "2,30 uur"
$ws.Range("D15").Value2 = 41365
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = "Steven V"
$ws.Range("G15").Value = "Solved"
$ws.Range("H15").Value = "APP"

# Move the active selection to G15, matching the saved view state
$ws.Range("G15").Select()
